# Generate Report for Handoff
# - Update the "Status" text from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it appears (Overview!E2/F2, zh-cn!C2,
#   de-de!C2).
# - Refresh the related timestamps (Overview!G2, de-de!H2 share one
#   value; zh-cn!H2 gets its own newer value).
# - Columns that displayed the long "Status" text are narrower now that
#   the new text is shorter, so re-fit them.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# Overview sheet: Status appears twice (zh-cn column, de-de column),
# followed by the shared "Latest HO Xliff Generate Date" timestamp.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-17 03:01:44"

# zh-cn detail sheet: Status column, then its own Latest Handoff Datetime.
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-17 03:01:39"

# de-de detail sheet: Status column, then Latest Handoff Datetime (same
# refreshed value as the Overview sheet).
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-17 03:01:44"

# The Status column got a lot shorter ("Handed back: in sync with
# en-US" -> "Ready for handoff"), so the report shrinks those columns
# to fit the new, narrower content instead of staying oversized.
$newStatusColumnWidth = 16.33
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
